$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.929.39'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.99%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.393.70'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.50%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.82%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.30'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.06%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.394.40'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.59%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("E9").Value = '  -0.73%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.57'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.70%  '

$ws.Range("E11").Value = '  -2.41%  '

$ws.Range("E12").Value = '  +1.35%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.970.81'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.53%  '

$ws.Range("E14").Value = '  +2.02%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.03'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.34%  '

$ws.Range("E16").Value = '  -1.43%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.388.32'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.56%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '60.996.93'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.02%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.14'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.86%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.84'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.30%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.96'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.86%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '383.03'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.65%  '

$ws.Range("E23").Value = '  -1.76%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.21%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.24%  '

$ws.Range("E26").Value = '  -5.07%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.531.62'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.59%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.179'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.66%  '

$ws.Range("E29").Value = '  -0.10%  '

$ws.Range("E30").Value = '  -3.10%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.02'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.72%  '

$ws.Range("E32").Value = '  -1.26%  '

$ws.Range("E33").Value = '  -2.51%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.47'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.01%  '

$ws.Range("E36").Value = '  -0.64%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '167.77'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.39%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.424.03'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.37%  '

$ws.Range("E39").Value = '  -2.41%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.48'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.03%  '

$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0772'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.51%  '

$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.45'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.85%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.782'
$ws.Range("D43").Style = "Normal"

$ws.Range("E44").Value = '  +0.04%  '

$ws.Range("E45").Value = '  -2.08%  '

$ws.Range("E46").Value = '  -3.95%  '

$ws.Range("E47").Value = '  -1.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.485.71'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.81%  '

$ws.Range("E49").Value = '  -2.03%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.02'
$ws.Range("D50").Style = "Normal"

$ws.Range("E51").Value = '  +1.02%  '
